# Ajout du TI du jour
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Luka Doncic
$ws.Range("F2").Value = 38.8
$ws.Range("G2").Value = 52.7
$ws.Range("H2").Value = 52.8
$ws.Range("N2").Value = 47
$ws.Range("O2").Value = 52
$ws.Range("P2").Value = 21
$ws.Range("Q2").Value = 51
$ws.Range("R2").Value = "-"
$ws.Range("U2").Value = -6.2

# Row 4 - Chet Holmgren
$ws.Range("I4").Value = 13
$ws.Range("M4").Value = 4

# Row 5 - Jonathan Kuminga
$ws.Range("I5").Value = 14
$ws.Range("L5").Value = 8

# Row 6 - D'Angelo Russell
# Clear the "Day-To-Day" status to an empty text cell (not a fully blank
# cell): force text-empty via the apostrophe text-prefix, then reset the
# style back to the sheet default so no stray formatting is left behind.
$ws.Range("C6").Value = "'"
$ws.Range("C6").Style = "Normal"
$ws.Range("I6").Value = 11
$ws.Range("L6").Value = 2

# Row 7 - Khris Middleton
$ws.Range("C7").Value = "Probable"

# Row 8 - Austin Reaves
$ws.Range("I8").Value = 12
$ws.Range("J8").Value = 1

# Row 9 - Stephen Curry
$ws.Range("I9").Value = 11
$ws.Range("J9").Value = 3

# Row 10 - Andrew Wiggins
$ws.Range("I10").Value = 10
$ws.Range("K10").Value = 3

# Row 11 - Gordon Hayward
$ws.Range("I11").Value = 13
$ws.Range("J11").Value = 13
